# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '35.335.99'
$ws.Cells.Item(2, 5).Value = '  -3.51%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.973.45'
$ws.Cells.Item(3, 5).Value = '  -4.83%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '240.18'
$ws.Cells.Item(5, 5).Value = '  -1.46%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.630'
$ws.Cells.Item(6, 5).Value = '  -5.07%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.18%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '55.91'
$ws.Cells.Item(8, 5).Value = '  +5.56%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '59.37'
$ws.Cells.Item(9, 5).Value = '  +1.06%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.356'
$ws.Cells.Item(10, 5).Value = '  -1.91%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -3.85%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.102'
$ws.Cells.Item(12, 5).Value = '  -6.28%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.890'
$ws.Cells.Item(13, 5).Value = '  -0.25%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '14.15'
$ws.Cells.Item(14, 5).Value = '  -3.86%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.265.37'
$ws.Cells.Item(15, 5).Value = '  -4.49%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -4.04%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '1.976.51'
$ws.Cells.Item(17, 5).Value = '  -4.73%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '17.03'
$ws.Cells.Item(18, 5).Value = '  +2.17%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '35.223.59'
$ws.Cells.Item(19, 5).Value = '  -3.69%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '69.74'
$ws.Cells.Item(20, 5).Value = '  -3.45%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.0₃0832'
$ws.Cells.Item(21, 5).Value = '  -3.78%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '231.18'
$ws.Cells.Item(22, 5).Value = '  -2.96%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.00'
$ws.Cells.Item(23, 5).Value = '  -6.59%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.03%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -5.71%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.23'
$ws.Cells.Item(26, 5).Value = '  +4.76%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '163.02'
$ws.Cells.Item(27, 5).Value = '  -1.41%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.01'
$ws.Cells.Item(28, 5).Value = '  -5.36%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '19.31'
$ws.Cells.Item(29, 5).Value = '  -6.20%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.118'
$ws.Cells.Item(30, 5).Value = '  -3.89%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.13'
$ws.Cells.Item(31, 5).Value = '  -2.20%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.75'
$ws.Cells.Item(32, 5).Value = '  -7.84%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -2.88%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0893'
$ws.Cells.Item(34, 5).Value = '  +10.30%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.24'
$ws.Cells.Item(35, 5).Value = '  -8.71%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.00%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.26'
$ws.Cells.Item(37, 5).Value = '  -5.41%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -2.99%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '4.81'
$ws.Cells.Item(39, 5).Value = '  -0.76%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -5.69%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.83'
$ws.Cells.Item(41, 5).Value = '  -1.60%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -5.01%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -6.28%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '90.00'
$ws.Cells.Item(44, 5).Value = '  -5.26%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0878'
$ws.Cells.Item(45, 5).Value = '  -7.55%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Maker'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.352.61'
$ws.Cells.Item(46, 5).Value = '  -2.86%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'FraxShare'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '7.37'
$ws.Cells.Item(47, 5).Value = '  -2.92%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '15.39'
$ws.Cells.Item(48, 5).Value = '  -1.44%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -0.79%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.25'
$ws.Cells.Item(50, 5).Value = '  -5.10%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '45.38'
$ws.Cells.Item(51, 5).Value = '  +0.24%  '
